$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.322.89"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").Value = "1.853.45"
$ws.Range("E3").Value = "  -3.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.50"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3949"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.74"
$ws.Range("E9").Value = "  -12.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07936"
$ws.Range("E10").Value = "  -5.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.010"
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.46"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "1.855.29"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.921"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.131"
$ws.Range("E15").Value = "  -4.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.12"
$ws.Range("E17").Value = "  -4.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001030"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06585"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  -4.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.470"
$ws.Range("E22").Value = "  -4.58%  "
$ws.Range("D23").Value = "27.328.38"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.297"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "2.064.62"
$ws.Range("E26").Value = "  -4.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.37"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.20"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.060"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.449"
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.50"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09416"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9474"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.438"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.587"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.260"
$ws.Range("E36").Value = "  -5.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06029"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02227"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.209"
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.025"
$ws.Range("E41").Value = "  -8.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5917"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1886"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.16"
$ws.Range("E44").Value = "  -8.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.284"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5614"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.06"
$ws.Range("E47").Value = "  -6.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.396"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.915"
$ws.Range("E49").Value = "  -5.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06762"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000308"
$ws.Range("E51").Value = "  +3.04%  "
